# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) column (F) across the sheets
# to the newly scraped values.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 2464
$ws.Range("F3").Value  = 754
$ws.Range("F4").Value  = 249
$ws.Range("F8").Value  = 913
$ws.Range("F14").Value = 74
$ws.Range("F16").Value = 1106
$ws.Range("F17").Value = 24416
$ws.Range("F18").Value = 2318
$ws.Range("F20").Value = 364
$ws.Range("F30").Value = 360
$ws.Range("F32").Value = 446

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value  = 272
$ws.Range("F8").Value  = 113
$ws.Range("F11").Value = 3642
$ws.Range("F13").Value = 156

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 798

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 2464
$ws.Range("F5").Value  = 798
$ws.Range("F6").Value  = 754
$ws.Range("F7").Value  = 249
$ws.Range("F14").Value = 272
$ws.Range("F16").Value = 913
$ws.Range("F21").Value = 74
$ws.Range("F23").Value = 1106
$ws.Range("F24").Value = 24416
$ws.Range("F28").Value = 156
$ws.Range("F30").Value = 2318
$ws.Range("F33").Value = 364
$ws.Range("F46").Value = 446
